$wb = $excel.ActiveWorkbook

# --- 1. Rename Sheet1 -> Diagram, add new "Square" sheet after it -------------
$diagram = $wb.Worksheets.Item(1)
$diagram.Name = "Diagram"
$square = $wb.Worksheets.Add($null, $diagram)
$square.Name = "Square"

# --- 2. Restructure the Diagram sheet row layout -------------------------------
# Insert 6 blank rows above the current data (old rows 1-15 become rows 7-21).
$diagram.Rows("1:6").Insert()

# The "Line" series (old rows 12-13) is now sitting at rows 18-19; move it up
# to rows 4-5 (it keeps its original formatting because Cut carries the style).
$diagram.Range("A18:C19").Cut($diagram.Range("A4"))
$diagram.Range("A18:C19").Clear()

# --- 3. New header row (row 1): left / top / size ------------------------------
$diagram.Range("A1").Style = "Normal"
$diagram.Range("B1").Value = "left"
$diagram.Range("C1").Value = "top"
$diagram.Range("D1").Value = "size"
$diagram.Range("B1:D1").Style = "Normal"

# --- 4. New box-definition rows (rows 2-3) --------------------------------------
$diagram.Range("A2").Value = "Square 1"
$diagram.Range("B2").Value = 10
$diagram.Range("C2").Value = 10
$diagram.Range("D2").Value = 10
$diagram.Range("A2:D2").Style = "Normal"

$diagram.Range("A3").Value = "Square 2"
$diagram.Range("B3").Value = 8
$diagram.Range("C3").Value = 7
$diagram.Range("D3").Value = 15
$diagram.Range("A3:D3").Style = "Normal"

# --- 5. Update the moved "Line" row values (rows 4-5) ---------------------------
$diagram.Range("B4").Value = 23
$diagram.Range("C4").Value = 7
$diagram.Range("B5").Value = 23
$diagram.Range("C5").Value = 7

# Row 6 is an empty spacer row between the "Line" block and the "Series" table.
$diagram.Range("A6:C6").Style = "Normal"

# --- 6. Square 1 point rows (rows 8-12) now reference the box definition -------
$diagram.Range("B8").Formula = "=Diagram!B2"
$diagram.Range("C8").Formula = "=Diagram!C2"
$diagram.Range("B9").Formula = "=B8"
$diagram.Range("C9").Formula = "=C8+Diagram!D2"
$diagram.Range("B10").Formula = "=B9+Diagram!D2"
$diagram.Range("C10").Formula = "=C9"
$diagram.Range("B11").Formula = "=B10"
$diagram.Range("C11").Formula = "=C10-Diagram!D2"
$diagram.Range("B12").Formula = "=B8"
$diagram.Range("C12").Formula = "=C8"
# B8 was filled, then dragged right across C8 - both ended up sharing B8's style.
$diagram.Range("C8").Style = $diagram.Range("B8").Style

# --- 7. Square 2 point rows (rows 13-17) now reference the box definition -----
$diagram.Range("B13").Formula = "=Diagram!B3"
$diagram.Range("C13").Formula = "=Diagram!C3"
$diagram.Range("B14").Formula = "=B13"
$diagram.Range("C14").Formula = "=C13+Diagram!D3"
$diagram.Range("B15").Formula = "=B14+Diagram!D3"
$diagram.Range("C15").Formula = "=C14"
$diagram.Range("B16").Formula = "=B15"
$diagram.Range("C16").Formula = "=C15-Diagram!D3"
$diagram.Range("B17").Formula = "=B13"
$diagram.Range("C17").Formula = "=C13"
# Same deal as B8/C8 above.
$diagram.Range("C13").Style = $diagram.Range("B13").Style

# --- 8. Middle 1 / Middle 2 rows (20-21) - formulas auto-shifted, reassert -----
$diagram.Range("B20").Formula = "=(MAX(B8:B12)+MIN(B8:B12))/2"
$diagram.Range("C20").Formula = "=(MAX(C8:C12)+MIN(C8:C12))/2"
$diagram.Range("B21").Formula = "=(MAX(B13:B17)+MIN(B13:B17))/2"
$diagram.Range("C21").Formula = "=(MAX(C13:C17)+MIN(C13:C17))/2"

# --- 9. Selection / view state ---------------------------------------------------
$diagram.Range("B2:D5").Select()
